# Daily attendance processing - 2025-12-16 15:56:40
#
# Normalizes the "Recorded By" (column G) entries on the active sheet:
# for every row whose value is a comma-separated list of recorders where
# the FIRST entry is "dnasr281@gmail.com" or "system" (case-insensitive),
# swap the first two entries (any further entries, e.g. "System", stay
# in place). Rows whose first entry is something else (e.g.
# "backup@backdoor.com", "admin@admin.com") or that only contain a single
# entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq "") {
        continue
    }

    $rawParts = $text.Split(",")
    if ($rawParts.Count -lt 2) {
        continue
    }

    $parts = @()
    for ($i = 0; $i -lt $rawParts.Count; $i++) {
        $parts += $rawParts[$i].Trim()
    }

    $firstLower = $parts[0].ToLower()
    if (($firstLower -eq "dnasr281@gmail.com") -or ($firstLower -eq "system")) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp

        $newText = $parts -join ", "
        $cell.Value = $newText
    }
}
